# issue #5: add legislator_id, name, date into dataframe
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- New header cells (H1:J1), styled like the existing header row ---
$hdr = $ws.Range("H1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- New data cells (row 2) ---
# Force the date column to stay plain text (otherwise Excel silently
# reinterprets "2011-11-21" as a date serial number).
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2011-11-21"

$ws.Range("I2").Value = "潘维剛"

# legislator_id is unknown for this record -> empty cell, but still
# present in the sheet (extends the used range to column J).
$ws.Range("J2").Borders.LineStyle = 0
